# Convert the "date" column (B) of the balance-sheet sheet from text dates
# (e.g. "1985-09-30") to real numeric Excel date serials, formatted as
# YYYY-MM-DD, matching the commit's "added ... datetime objects for all
# dataframes" change.

function ConvertTo-ExcelSerial($y, $m, $d) {
    # Howard Hinnant's civil_from_days algorithm (run in reverse), days
    # since 1970-01-01, then re-based onto Excel's serial-date epoch
    # (1899-12-30 == serial 0). Pure integer arithmetic -- avoids relying
    # on a [DateTime] class.
    $yy = $y
    if ($m -le 2) {
        $yy = $yy - 1
    }
    if ($yy -ge 0) {
        $era = [Math]::Floor($yy / 400)
    } else {
        $era = [Math]::Floor(($yy - 399) / 400)
    }
    $yoe = $yy - $era * 400
    $mp = ($m + 9) % 12
    $doy = [Math]::Floor((153 * $mp + 2) / 5) + $d - 1
    $doe = $yoe * 365 + [Math]::Floor($yoe / 4) - [Math]::Floor($yoe / 100) + $doy
    $daysSinceUnixEpoch = $era * 146097 + $doe - 719468
    return $daysSinceUnixEpoch + 25569
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $text = $cell.Text
    if ([string]::IsNullOrEmpty($text)) {
        continue
    }
    $parts = $text.Split("-")
    $y = [int]$parts[0]
    $m = [int]$parts[1]
    $d = [int]$parts[2]
    $serial = ConvertTo-ExcelSerial $y $m $d
    $cell.Value = $serial
    $cell.NumberFormat = "YYYY-MM-DD"
}

Write-Host "Converted date column B2:B$lastRow to numeric dates"
